# config file change for docker work
$wb = $excel.ActiveWorkbook

# RUNMANAGER sheet: flip the "execute" flags for rows 3-6 from "no" to "yes"
$wsRunManager = $wb.Worksheets.Item("RUNMANAGER")
$wsRunManager.Range("C3").Value = "yes"
$wsRunManager.Range("C4").Value = "yes"
$wsRunManager.Range("C5").Value = "yes"
$wsRunManager.Range("C6").Value = "yes"
$wsRunManager.Range("C2:C6").Select()

# DATA sheet: flip the "execute" flags for rows 3-6 from "no" to "yes"
$wsData = $wb.Worksheets.Item("DATA")
$wsData.Range("B3").Value = "yes"
$wsData.Range("B4").Value = "yes"
$wsData.Range("B5").Value = "yes"
$wsData.Range("B6").Value = "yes"
$wsData.Range("B2:B6").Select()

# DATA becomes the active sheet/tab
$wsData.Activate()
